# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (D) / "Correspond Handback DateTime" (G)
# timestamps on row 5 (the 67741593-... handback file) for both the zh-cn and
# de-de language sheets, advancing them to reflect the freshly regenerated
# handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-25 09:06:38"
$wsZhCn.Range("G5").Value = "2016-02-25 09:07:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-25 09:06:48"
$wsDeDe.Range("G5").Value = "2016-02-25 09:07:45"
